$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (PR) - meta_avg (H) and meta_std (I) updated
$ws.Range("H17").Value = 1266.712068611907
$ws.Range("I17").Value = 1414.634209378984

# Row 18 - meta_avg (H) and meta_std (I) updated
$ws.Range("H18").Value = 1211.184666308004
$ws.Range("I18").Value = 1225.698890613596

# Row 21 - meta_avg (H) and meta_std (I) updated
$ws.Range("H21").Value = 826.0106721316131
$ws.Range("I21").Value = 949.1475679611532

# Row 22 - meta_avg (H), meta_std (I), meta_min (J) updated
$ws.Range("H22").Value = 1714.950834610086
$ws.Range("I22").Value = 2425.306729108642
$ws.Range("J22").Value = 0

# Row 24 - meta_avg (H) and meta_std (I) updated
$ws.Range("H24").Value = 1251.886692580295
$ws.Range("I24").Value = 3208.981321493234
